$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date serial for rows 2 through 28 (45418 -> 45419)
for ($r = 2; $r -le 28; $r++) {
  $ws.Cells.Item($r, 3).Value = 45419
}

# Remove the last data row (row 29, "A 17575-2024") entirely
$ws.Rows.Item(29).Delete()

# Row 28 no longer needs an explicit custom height; reset it to the default
$ws.Rows.Item(28).AutoFit()
